$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record data for each data row (rows 2-57)
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD = column 30
    $ws.Cells.Item($r, 31).Value = 70  # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0   # AF = column 32
}
